$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'62.404.94"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -2.49%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.004.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -2.99%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.09%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'583.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.46%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'146.53"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -6.40%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -3.14%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'3.003.88"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.09%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.149"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -5.66%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'5.71"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -3.53%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.84%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.59%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'34.68"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -6.13%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.123"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.16%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.497.65"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.08%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'7.02"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = "'62.399.82"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -2.37%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'3.002.96"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.08%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'458.86"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -4.51%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'13.88"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -4.02%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -4.76%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'7.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -2.54%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'80.02"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.70%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -8.63%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'12.23"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -5.44%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.14%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'10.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -6.88%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.04%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'7.17"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -5.33%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -2.68%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -4.98%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.20%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -5.26%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -3.80%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.0₃0790"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -6.20%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -4.95%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -6.52%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'50.00"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -2.02%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = "'dogwifhat"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'2.94"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -10.55%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'Cosmos"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'8.92"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -3.73%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'409.14"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -7.59%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.83%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.276"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -4.93%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'39.13"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -2.61%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'2.770.89"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.26%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.0353"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -3.23%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'127.45"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.99%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.05%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.94%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'23.66"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -9.41%  "
$ws.Range('E51').Style = 'Normal'
